$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.543.22"
$ws.Range("E2").Value = "  +3.59%  "

$ws.Range("D3").Value = "3.503.67"
$ws.Range("E3").Value = "  +2.28%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.09%  "

$ws.Range("D8").Value = "3.498.88"
$ws.Range("E8").Value = "  +2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.127"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.33%  "

$ws.Range("D13").Value = "4.106.55"
$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.73%  "

$ws.Range("E16").Value = "  +1.88%  "

$ws.Range("D17").Value = "66.565.17"

$ws.Range("D18").Value = "3.505.03"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.38%  "

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.537"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.182"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.55%  "

$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.07%  "

$ws.Range("E32").Value = "  +2.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.30%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +7.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.884"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0749"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.83%  "

$ws.Range("D45").Value = "2.797.57"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("E47").Value = "  +1.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "351.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.57%  "

$ws.Range("E50").Value = "  +3.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.75%  "
